# Update cryptos list: price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.058.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").Value = "'1.668.30"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'216.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").Value = "'0.5110"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.2672"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").Value = "'0.06398"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").Value = "'21.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "'0.07450"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "'1.692.43"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "'4.518"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "'0.5814"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "'0.000008509"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("D16").Value = "'64.15"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("D17").Value = "'25.883.97"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("D18").Value = "'4.929"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "'10.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").Value = "'189.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").Value = "'6.197"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'145.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").Value = "'7.622"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("D26").Value = "'0.1222"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.37%  "

$ws.Range("D27").Value = "'15.65"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("D28").Value = "'0.06701"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +16.53%  "

$ws.Range("D29").Value = "'1.332"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.29%  "

$ws.Range("D30").Value = "'1.312"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").Value = "'3.555"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.87%  "

$ws.Range("D32").Value = "'3.516"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "'1.662"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").Value = "'1.019"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").Value = "'0.6171"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.51%  "

$ws.Range("D36").Value = "'2.368"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").Value = "'2.682"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("D38").Value = "'6.346"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.63%  "

$ws.Range("D39").Value = "'1.097.34"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "'0.8693"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").Value = "'101.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.61%  "

$ws.Range("D44").Value = "'1.814.41"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").Value = "'0.00000000116"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46").Value = "'56.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "'8.103"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'0.05232"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").Value = "'6.004"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.04%  "
